# Last Update 17-09-2018 14:02:18.62
# Fill in the remaining attendance hours (M,N columns) on the Attendance sheet
# and move the active-tab focus back to the Attendance sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attendance")
$labWs = $wb.Worksheets.Item("Lab Attendance")

# --- Header hour totals for the new columns (row 4: total students, row 5/6: breakdown) ---
$ws.Cells.Item(4, 13).Value = 9     # M4
$ws.Cells.Item(4, 14).Value = 9     # N4
$ws.Cells.Item(5, 13).Value = 15    # M5
$ws.Cells.Item(5, 14).Value = 17    # N5
$ws.Cells.Item(6, 13).Value = 7     # M6
$ws.Cells.Item(6, 14).Value = 5     # N6

# --- Daily attendance (P = present, A = absent) for the two newly-recorded hours ---
$attendance = @{
    7  = @("A","P")
    8  = @("A","P")
    9  = @("P","P")
    10 = @("P","P")
    11 = @("P","P")
    12 = @("A","P")
    13 = @("P","P")
    14 = @("P","P")
    15 = @("P","A")
    16 = @("A","P")
    17 = @("P","A")
    18 = @("P","P")
    19 = @("P","P")
    20 = @("P","P")
    21 = @("A","P")
    22 = @("A","P")
    23 = @("A","A")
    24 = @("P","P")
    25 = @("P","P")
    26 = @("P","P")
    27 = @("P","P")
    28 = @("P","P")
    29 = @("P","P")
    30 = @("A","A")
    31 = @("P","P")
    32 = @("A","P")
    33 = @("P","P")
    34 = @("P","P")
    35 = @("P","P")
    36 = @("P","P")
    37 = @("P","P")
    38 = @("P","P")
    39 = @("P","P")
    40 = @("A","P")
    41 = @("P","P")
    42 = @("P","P")
    43 = @("P","P")
    44 = @("A","P")
    45 = @("A","P")
    46 = @("P","P")
    47 = @("P","P")
    48 = @("P","P")
    49 = @("P","P")
    50 = @("P","P")
    51 = @("P","P")
    52 = @("P","P")
    53 = @("P","P")
    54 = @("P","P")
    55 = @("P","P")
}

foreach ($row in $attendance.Keys) {
    $vals = $attendance[$row]
    $ws.Cells.Item($row, 13).Value = $vals[0]
    $ws.Cells.Item($row, 14).Value = $vals[1]
}

# --- Restore the Attendance sheet as the active tab / selection, matching where the
#     author left off editing (scrolled down near the bottom of the table). ---
$ws.Activate()
$ws.Range("N56").Select()
$excel.ActiveWindow.ScrollRow = 32
$excel.ActiveWindow.ScrollColumn = 1

# The "Lab Attendance" sheet's own lingering selection moves off of its previous
# cell now that it is no longer the active tab.
$labWs.Range("A25").Select()

$ws.Activate()
